$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "45.376.73"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "  +5.60%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "2.366.86"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = "  +2.88%  "
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "  +0.06%  "
$c.Style = "Normal"

$c = $ws.Range("B5")
$c.Value = "BNB"
$c.Style = "Normal"

$c = $ws.Range("C5")
$c.Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'313.86"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = "  +1.27%  "
$c.Style = "Normal"

$c = $ws.Range("B6")
$c.Value = "Solana"
$c.Style = "Normal"

$c = $ws.Range("C6")
$c.Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'110.04"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "  +4.90%  "
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "  +0.90%  "
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "  +0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.614"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "  +1.64%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'40.93"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "  +3.45%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.0917"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "  +1.01%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'8.51"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "  +2.81%  "
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "  +1.75%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'0.980"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "  -0.57%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "2.730.28"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "  +2.97%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'15.46"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "  +1.34%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "2.366.66"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "  +3.41%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "45.333.85"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "  +5.88%  "
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "  -0.38%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.0000106"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "  +1.65%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'13.07"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "  -5.20%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'73.60"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "  +0.29%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'3.47"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "  +0.59%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'261.40"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "  -2.82%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'2.29"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "  +1.73%  "
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "  -0.61%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'11.11"
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "  +1.62%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'7.35"
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "  -5.99%  "
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "  +2.25%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'22.47"
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.Value = "  +1.22%  "
$c.Style = "Normal"

$c = $ws.Range("B31")
$c.Value = "InjectiveProtocol"
$c.Style = "Normal"

$c = $ws.Range("C31")
$c.Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'37.96"
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = "  -0.02%  "
$c.Style = "Normal"

$c = $ws.Range("B32")
$c.Value = "Hedera"
$c.Style = "Normal"

$c = $ws.Range("C32")
$c.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'0.0960"
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "  +11.07%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'169.68"
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "  +2.25%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'2.95"
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = "  +4.66%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'0.130"
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = "  -0.04%  "
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "  +2.80%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'4.80"
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.Value = "  +3.88%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'3.98"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = "  +10.41%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'2.97"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "  +6.91%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.0355"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "  -0.68%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'1.72"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "  +10.32%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'102.12"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "  -5.06%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.234"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "  +2.55%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'13.23"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "  +8.24%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'69.92"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "  -1.48%  "
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "  -0.33%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'81.05"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "  +6.86%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'112.65"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "  +1.18%  "
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'9.31"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "  +5.21%  "
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "  +7.16%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "1.642.49"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "  -3.27%  "
$c.Style = "Normal"

Write-Host "Applied 99 cell updates"